{"js": "// Update the worksheet: refresh the header date and every three-digit\n// \u00f7 one-digit practice problem (new operands/quotients/remainders),\n// per the \"c8c62b6\" output regeneration.\nconst replacements = [\n  [\"2026-02-28 Saturday\", \"2026-03-01 Sunday\"],\n  [\"778\u00f79=86, 4\", \"440\u00f74=110, 0\"],\n  [\"413\u00f72=206, 1\", \"379\u00f72=189, 1\"],\n  [\"754\u00f73=251, 1\", \"817\u00f75=163, 2\"],\n  [\"843\u00f78=105, 3\", \"299\u00f76=49, 5\"],\n  [\"389\u00f73=129, 2\", \"983\u00f79=109, 2\"],\n  [\"684\u00f77=97, 5\", \"441\u00f77=63, 0\"],\n  [\"259\u00f79=28, 7\", \"471\u00f75=94, 1\"],\n  [\"825\u00f75=165, 0\", \"216\u00f75=43, 1\"],\n  [\"905\u00f78=113, 1\", \"878\u00f75=175, 3\"],\n  [\"841\u00f78=105, 1\", \"968\u00f75=193, 3\"],\n  [\"199\u00f78=24, 7\", \"555\u00f78=69, 3\"],\n  [\"319\u00f74=79, 3\", \"621\u00f73=207, 0\"],\n  [\"712\u00f77=101, 5\", \"123\u00f73=41, 0\"],\n  [\"991\u00f76=165, 1\", \"894\u00f77=127, 5\"],\n  [\"372\u00f73=124, 0\", \"546\u00f72=273, 0\"],\n  [\"188\u00f78=23, 4\", \"738\u00f72=369, 0\"],\n  [\"127\u00f72=63, 1\", \"960\u00f76=160, 0\"],\n  [\"543\u00f79=60, 3\", \"557\u00f75=111, 2\"],\n  [\"214\u00f74=53, 2\", \"711\u00f75=142, 1\"],\n  [\"421\u00f76=70, 1\", \"688\u00f78=86, 0\"],\n  [\"395\u00f78=49, 3\", \"317\u00f73=105, 2\"],\n  [\"233\u00f73=77, 2\", \"321\u00f73=107, 0\"],\n  [\"743\u00f73=247, 2\", \"967\u00f73=322, 1\"],\n  [\"138\u00f72=69, 0\", \"618\u00f75=123, 3\"],\n  [\"145\u00f76=24, 1\", \"455\u00f78=56, 7\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet: refresh the header date and every three-digit\n# \u00f7 one-digit practice problem (new operands/quotients/remainders),\n# per the \"c8c62b6\" output regeneration.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('2026-02-28 Saturday', '2026-03-01 Sunday'),\n    @('778\u00f79=86, 4', '440\u00f74=110, 0'),\n    @('413\u00f72=206, 1', '379\u00f72=189, 1'),\n    @('754\u00f73=251, 1', '817\u00f75=163, 2'),\n    @('843\u00f78=105, 3', '299\u00f76=49, 5'),\n    @('389\u00f73=129, 2', '983\u00f79=109, 2'),\n    @('684\u00f77=97, 5', '441\u00f77=63, 0'),\n    @('259\u00f79=28, 7', '471\u00f75=94, 1'),\n    @('825\u00f75=165, 0', '216\u00f75=43, 1'),\n    @('905\u00f78=113, 1', '878\u00f75=175, 3'),\n    @('841\u00f78=105, 1', '968\u00f75=193, 3'),\n    @('199\u00f78=24, 7', '555\u00f78=69, 3'),\n    @('319\u00f74=79, 3', '621\u00f73=207, 0'),\n    @('712\u00f77=101, 5', '123\u00f73=41, 0'),\n    @('991\u00f76=165, 1', '894\u00f77=127, 5'),\n    @('372\u00f73=124, 0', '546\u00f72=273, 0'),\n    @('188\u00f78=23, 4', '738\u00f72=369, 0'),\n    @('127\u00f72=63, 1', '960\u00f76=160, 0'),\n    @('543\u00f79=60, 3', '557\u00f75=111, 2'),\n    @('214\u00f74=53, 2', '711\u00f75=142, 1'),\n    @('421\u00f76=70, 1', '688\u00f78=86, 0'),\n    @('395\u00f78=49, 3', '317\u00f73=105, 2'),\n    @('233\u00f73=77, 2', '321\u00f73=107, 0'),\n    @('743\u00f73=247, 2', '967\u00f73=322, 1'),\n    @('138\u00f72=69, 0', '618\u00f75=123, 3'),\n    @('145\u00f76=24, 1', '455\u00f78=56, 7'),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    # wdFindContinue=1, wdReplaceAll=2\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"Text not found: $oldText\"\n    }\n}\n"}
